$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values with new TPM-derived figures
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.370913
$ws.Range("H2").Value = 1.112739
$ws.Range("M2").Value = 0.3045636666666667
$ws.Range("N2").Value = 0.913691
$ws.Range("Q2").Value = 0.1129666232943333
$ws.Range("R2").Value = 1.016699609649
